$wb = $excel.ActiveWorkbook

# --- Overview sheet: mirror the new status text (shared string reused by B2/C2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

# --- Per-language sheets ---
$langSheets = @{
    "zh-cn" = "https://github.com/OpenLocalizationTest/oltest/blob/c4476c7ea412af34f42a1168f8c01738edd3f7e8/e2e/3d029bba-3824-48aa-ba83-1438ac837909.md"
    "de-de" = "https://github.com/OpenLocalizationTest/oltest/blob/c4476c7ea412af34f42a1168f8c01738edd3f7e8/e2e/3d029bba-3824-48aa-ba83-1438ac837909.md"
}
$configTarget = "https://github.com/OpenLocalizationTest/oltest/blob/c4476c7ea412af34f42a1168f8c01738edd3f7e8/.localization-config"

foreach ($sheetName in @("zh-cn","de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff File (C2) is no longer available -> drop the cell + its hyperlink.
    # This shim's Hyperlinks.Delete() clears every hyperlink on the sheet, so capture
    # the two links we need to keep, wipe them all, then recreate those two.
    $mdDisplay = $ws.Range("A2").Value2
    $cfgDisplay = $ws.Range("A3").Value2

    $ws.Hyperlinks.Delete()

    $ws.Range("C2").Clear()

    $ws.Hyperlinks.Add($ws.Range("A2"), $langSheets[$sheetName], "", "", $mdDisplay)
    $ws.Hyperlinks.Add($ws.Range("A3"), $configTarget, "", "", $cfgDisplay)

    # Latest Handoff Datetime (D2) now mirrors the "unset" sentinel datetime
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason (G2) stays the sentinel datetime value
    $ws.Range("G2").Value = "0001-01-01 00:00:00"

    # Dependency From (H2) -> "Ignored"
    $ws.Range("H2").Value = "Ignored"

    # Row 3 values are unchanged from before (still the sentinel datetime / Ignored)
    $ws.Range("D3").Value = "0001-01-01 00:00:00"
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Ignored"
}
